$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "26.474.37"
$ws.Range("E2").Value = "  +1.50%  "

# Row 3
$ws.Range("D3").Value = "1.677.31"
$ws.Range("E3").Value = "  +2.28%  "

# Row 4
$ws.Range("E4").Value = "  -0.04%  "

# Row 5
$ws.Range("D5").Value = "'219.92"
$ws.Range("E5").Value = "  +2.55%  "

# Row 6
$ws.Range("D6").Value = "'0.5317"
$ws.Range("E6").Value = "  +2.02%  "

# Row 7
$ws.Range("D7").Value = "'1.001"
$ws.Range("E7").Value = "  -0.07%  "

# Row 8
$ws.Range("D8").Value = "'0.2698"
$ws.Range("E8").Value = "  +3.70%  "

# Row 9
$ws.Range("D9").Value = "'0.06405"
$ws.Range("E9").Value = "  +1.54%  "

# Row 10
$ws.Range("D10").Value = "'21.85"
$ws.Range("E10").Value = "  +5.61%  "

# Row 11
$ws.Range("D11").Value = "'0.07802"
$ws.Range("E11").Value = "  +1.54%  "

# Row 12
$ws.Range("B12").Value = "Polkadot"
$ws.Range("C12").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D12").Value = "'4.514"
$ws.Range("E12").Value = "  +2.32%  "

# Row 13
$ws.Range("B13").Value = "WrappedEther"
$ws.Range("C13").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D13").Value = "1.675.34"
$ws.Range("E13").Value = "  +3.25%  "

# Row 14
$ws.Range("D14").Value = "'0.5594"
$ws.Range("E14").Value = "  +0.66%  "

# Row 15
$ws.Range("D15").Value = "0.0₅8348"
$ws.Range("E15").Value = "  +1.98%  "

# Row 16
$ws.Range("D16").Value = "'65.70"
$ws.Range("E16").Value = "  +0.99%  "

# Row 17
$ws.Range("D17").Value = "26.498.96"
$ws.Range("E17").Value = "  +1.59%  "

# Row 18
$ws.Range("E18").Value = "  -0.10%  "

# Row 19
$ws.Range("D19").Value = "'4.792"
$ws.Range("E19").Value = "  +1.53%  "

# Row 20
$ws.Range("D20").Value = "'193.23"
$ws.Range("E20").Value = "  +2.17%  "

# Row 21
$ws.Range("D21").Value = "'10.31"
$ws.Range("E21").Value = "  +1.09%  "

# Row 22
$ws.Range("D22").Value = "'6.332"
$ws.Range("E22").Value = "  +2.63%  "

# Row 23
$ws.Range("E23").Value = "  -0.06%  "

# Row 24
$ws.Range("D24").Value = "'0.1275"
$ws.Range("E24").Value = "  +5.97%  "

# Row 25
$ws.Range("D25").Value = "'140.14"
$ws.Range("E25").Value = "  -3.72%  "

# Row 26
$ws.Range("D26").Value = "'7.416"
$ws.Range("E26").Value = "  -0.01%  "

# Row 27
$ws.Range("D27").Value = "'16.30"
$ws.Range("E27").Value = "  +2.92%  "

# Row 28
$ws.Range("D28").Value = "'1.445"
$ws.Range("E28").Value = "  +3.68%  "

# Row 29
$ws.Range("D29").Value = "'0.06267"
$ws.Range("E29").Value = "  +6.49%  "

# Row 30
$ws.Range("D30").Value = "'1.291"
$ws.Range("E30").Value = "  +2.87%  "

# Row 31
$ws.Range("D31").Value = "'3.614"
$ws.Range("E31").Value = "  +5.04%  "

# Row 32
$ws.Range("D32").Value = "'3.458"
$ws.Range("E32").Value = "  +1.59%  "

# Row 33
$ws.Range("D33").Value = "'1.697"
$ws.Range("E33").Value = "  +2.95%  "

# Row 34
$ws.Range("D34").Value = "'1.014"
$ws.Range("E34").Value = "  +3.11%  "

# Row 35
$ws.Range("D35").Value = "'0.6193"
$ws.Range("E35").Value = "  +9.65%  "

# Row 36
$ws.Range("D36").Value = "'2.420"
$ws.Range("E36").Value = "  +1.15%  "

# Row 37
$ws.Range("E37").Value = "  +0.86%  "

# Row 38
$ws.Range("D38").Value = "'6.164"
$ws.Range("E38").Value = "  +8.03%  "

# Row 39
$ws.Range("D39").Value = "'0.01635"
$ws.Range("E39").Value = "  +1.21%  "

# Row 40
$ws.Range("D40").Value = "1.097.19"
$ws.Range("E40").Value = "  +6.86%  "

# Row 41
$ws.Range("D41").Value = "'0.8623"
$ws.Range("E41").Value = "  +0.99%  "

# Row 42
$ws.Range("D42").Value = "'1.0000"
$ws.Range("E42").Value = "  -0.13%  "

# Row 43
$ws.Range("D43").Value = "'100.67"
$ws.Range("E43").Value = "  +0.57%  "

# Row 44
$ws.Range("E44").Value = "  +1.69%  "

# Row 45
$ws.Range("B45").Value = "Aave"
$ws.Range("C45").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D45").Value = "'58.75"
$ws.Range("E45").Value = "  +5.32%  "

# Row 46
$ws.Range("B46").Value = "BabyDogeCoin"
$ws.Range("C46").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("D46").Value = "0.0₈109"
$ws.Range("E46").Value = "  +2.10%  "

# Row 47
$ws.Range("D47").Value = "'8.171"
$ws.Range("E47").Value = "  +1.42%  "

# Row 48
$ws.Range("D48").Value = "'1.006"
$ws.Range("E48").Value = "  +0.31%  "

# Row 49
$ws.Range("B49").Value = "Cronos"
$ws.Range("C49").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D49").Value = "'0.05194"
$ws.Range("E49").Value = "  +0.86%  "

# Row 50
$ws.Range("B50").Value = "RenderToken"
$ws.Range("C50").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D50").Value = "'1.482"
$ws.Range("E50").Value = "  +7.22%  "

# Row 51
$ws.Range("D51").Value = "'6.037"
$ws.Range("E51").Value = "  +2.35%  "

